# Weekly data refresh: insert the newest observation as a new row right
# above the existing row 21, shifting the rest of the "Vega Modelo de
# Temuco" / Guayaba block down by one row (old row 21 -> 22, ... old row
# 33 -> 34). This matches the sheet growing from A1:T33 to A1:T34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21 (pushes 21..33 down to 22..34).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the latest market observation.
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 45086
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108001
$ws.Range("J21").Value = "Guayaba"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 30
$ws.Range("N21").Value = 2600
$ws.Range("O21").Value = 2600
$ws.Range("P21").Value = 2600
$ws.Range("Q21").Value = "$/kilo"
$ws.Range("R21").Value = "Región de Arica y Parinacota"
$ws.Range("S21").Value = 2600
$ws.Range("T21").Value = 1
